# Clean up the weekly-report paragraphs: the original text had been
# spell-checked by Word (w:proofErr spellStart/spellEnd/gramStart/gramEnd
# markers splitting almost every word into its own run). This pass
# retypes each affected paragraph's text as a single run, and relocates
# the trailing "_GoBack" bookmark into its own empty paragraph.
#
# Note: a plain Range.Text assignment is a no-op when the new string
# equals the existing text, so the proofErr-riddled runs are collapsed
# with Find.Execute (a real find/replace) instead, scoped to each
# paragraph's own Range so nothing outside it is touched.

$d = $word.ActiveDocument

function Set-ParagraphText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $null = $r.Find.Execute($text, $true, $false, $false, $false, $false, `
        $true, 1, $false, $text, 2)
}

# Paragraph 2: "Đại Anh Dũng : "
Set-ParagraphText 2 "Đại Anh Dũng : "

# Paragraph 3: "_Dựng layout web bán giày"
Set-ParagraphText 3 "_Dựng layout web bán giày"

# Paragraph 4: "_Dựng page login, signup, forgot password" followed by an
# already-separate run ", home, detail, cart" that must stay its own run.
$p4 = $d.Paragraphs.Item(4)
$head = $p4.Range
$null = $head.Find.Execute("_Dựng page login, signup, forgot password", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "_Dựng page login, signup, forgot password", 2)
# Retyping the head merged every run of the paragraph into one (including
# the previously-separate trailing run). Re-split the trailing run back
# off by nudging a character property on it and back to what it was.
$p4 = $d.Paragraphs.Item(4)
$headLen = "_Dựng page login, signup, forgot password".Length
$tailStart = $p4.Range.Start + $headLen
$tail = $d.Range($tailStart, $p4.Range.End - 1)
$tail.Font.Bold = $true
$tail.Font.Bold = $false

# Paragraph 5: "_Get Api login, signup, forgot password"
Set-ParagraphText 5 "_Get Api login, signup, forgot password"

# Paragraph 6: "_Fix api endpoint"
Set-ParagraphText 6 "_Fix api endpoint"

# Paragraph 7: "_Tạo api upload get image"
Set-ParagraphText 7 "_Tạo api upload get image"

# Paragraph 10: "Đại Anh Dũng"
Set-ParagraphText 10 "Đại Anh Dũng"

# Paragraph 11: "_Viết, get Api home slide, detail" — also push the
# "_GoBack" bookmark that trails it onto a new paragraph of its own by
# replacing the text with itself plus a paragraph mark.
$found = $d.Content.Find.Execute(
    "_Viết, get Api home slide, detail", $true, $false, $false, $false, `
    $false, $true, 1, $false, "_Viết, get Api home slide, detail^p", 2)
